$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42636.589062500003
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9917.16
$ws.Range("D5").Value = 9948
$ws.Range("E5").Value = 19.29
$ws.Range("F5").Value = 19.41

$ws.Range("G5").Value = $true
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"

$ws.Range("H5").Value = 0.62
$ws.Range("I5").Value = $false
